# Auto-generated Excel COM-interop edit script.
# Updates cached price/profit figures on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match the latest scheduled-runner snapshot. All target cells hold literal numbers
# (no formulas), so each change is a direct value write; cells that the new snapshot
# omits are cleared, and cells newly present in the snapshot are written fresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: H33,I33,K33,M33
$ws.Range("H33").Value = 95.48148
$ws.Range("I33").Value = 108.1875
$ws.Range("K33").Value = 108.1875
$ws.Range("M33").Value = 120.8125
# Row 58: H58,I58,K58,M58
$ws.Range("H58").Value = 2198.2
$ws.Range("I58").Value = 330.33334
$ws.Range("K58").Value = 991.0000200000001
$ws.Range("M58").Value = -841.0000200000001
# Row 64: H64,I64,K64,M64
$ws.Range("H64").Value = 4407
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
# Row 67: H67,I67,K67,M67
$ws.Range("H67").Value = 4407
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
# Row 69: H69,I69,J69,K69,L69,M69,N69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
# Row 72: H72,I72,J72,K72,L72,M72,N72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
# Row 99: H99,I99,K99,M99
$ws.Range("H99").Value = 542.8182
$ws.Range("I99").Value = 453.6
$ws.Range("K99").Value = 1360.8
$ws.Range("M99").Value = 137.1999999999998
# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 2962.875
$ws.Range("I107").Value = 2375.5
$ws.Range("K107").Value = 2375.5
$ws.Range("M107").Value = -455.5
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1450
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 4350
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -9410
# Row 137: H137,I137,K137,M137
$ws.Range("H137").Value = 2814.5
$ws.Range("I137").Value = 2814.5
$ws.Range("K137").Value = 8443.5
$ws.Range("M137").Value = -5893.5
# Row 138: H138,J138,L138,N138
$ws.Range("H138").Value = 2581.4614
$ws.Range("J138").Value = 3933
$ws.Range("L138").Value = 11799
$ws.Range("N138").Value = -22079

$ws = $wb.Worksheets.Item("ARM")
# Row 46: H46,J46,L46,N46
$ws.Range("H46").Value = 3999.5
$ws.Range("J46").Value = 3999.5
$ws.Range("L46").Value = 3999.5
$ws.Range("N46").Value = -4637.5
# Row 63: H63,I63,K63,M63
$ws.Range("H63").Value = 1995
$ws.Range("I63").Value = 1995
$ws.Range("K63").Value = 1995
$ws.Range("M63").Value = -1309
# Row 66: H66,I66,K66,M66
$ws.Range("H66").Value = 1995
$ws.Range("I66").Value = 1995
$ws.Range("K66").Value = 9975
$ws.Range("M66").Value = -6543
# Row 74: H74,I74,K74,M74
$ws.Range("H74").Value = 1466
$ws.Range("I74").Value = 1466
$ws.Range("K74").Value = 1466
$ws.Range("M74").Value = -592
# Row 77: H77,I77,K77,M77
$ws.Range("H77").Value = 1466
$ws.Range("I77").Value = 1466
$ws.Range("K77").Value = 7330
$ws.Range("M77").Value = -2962
# Row 130: H130,J130,L130,N130
$ws.Range("H130").Value = 65428
$ws.Range("J130").Value = 65428
$ws.Range("L130").Value = 65428
$ws.Range("N130").Value = -75468

$ws = $wb.Worksheets.Item("BSM")
# Row 22: H22,J22,L22,N22
$ws.Range("H22").Value = 174
$ws.Range("J22").Value = 169
$ws.Range("L22").Value = 169
$ws.Range("N22").Value = -515
# Row 62: H62,J62,L62,N62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65: H65,J65,L65,N65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 82: H82,J82,L82,N82
$ws.Range("H82").Value = 23190.646
$ws.Range("J82").Value = 29998.334
$ws.Range("L82").Value = 29998.334
$ws.Range("N82").Value = -30764.334
# Row 85: H85,J85,L85,N85
$ws.Range("H85").Value = 23190.646
$ws.Range("J85").Value = 29998.334
$ws.Range("L85").Value = 29998.334
$ws.Range("N85").Value = -32650.334

$ws = $wb.Worksheets.Item("CRP")
# Row 7: H7,I7,K7,M7
$ws.Range("H7").Value = 259.0909
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 100
$ws.Range("M7").Value = 13
# Row 58: H58,I58,J58,K58,L58,M58,N58
$ws.Range("H58").Value = 1613.375
$ws.Range("I58").Value = 1623.4546
$ws.Range("J58").Value = 1604.8462
$ws.Range("K58").Value = 1623.4546
$ws.Range("L58").Value = 1604.8462
$ws.Range("M58").Value = -1420.4546
$ws.Range("N58").Value = -2010.8462
# Row 105: H105,I105,K105,M105
$ws.Range("H105").Value = 865.6667
$ws.Range("I105").Value = 1065
$ws.Range("K105").Value = 1065
$ws.Range("M105").Value = 682
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 2966.818
$ws.Range("I132").Value = 2965.5
$ws.Range("K132").Value = 8896.5
$ws.Range("M132").Value = -6366.5
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 3294.889
$ws.Range("I134").Value = 3145.75
$ws.Range("K134").Value = 9437.25
$ws.Range("M134").Value = -6902.25
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 1613.375
$ws.Range("I136").Value = 1623.4546
$ws.Range("J136").Value = 1604.8462
$ws.Range("K136").Value = 4870.3638
$ws.Range("L136").Value = 4814.5386
$ws.Range("M136").Value = -2320.3638
$ws.Range("N136").Value = -9914.5386

$ws = $wb.Worksheets.Item("CUL")
# Row 4: H4,J4,L4,N4
$ws.Range("H4").Value = 1498.1904
$ws.Range("J4").Value = 2040
$ws.Range("L4").Value = 6120
$ws.Range("N4").Value = -6344
# Row 109: H109,I109,K109,M109
$ws.Range("H109").Value = 698.3333
$ws.Range("I109").Value = 698.3333
$ws.Range("K109").Value = 2094.9999
$ws.Range("M109").Value = -1054.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 2: H2,I2,J2,K2,L2,M2,N2
$ws.Range("H2").Value = 480.9565
$ws.Range("I2").Value = 551.375
$ws.Range("J2").Value = 320
$ws.Range("K2").Value = 551.375
$ws.Range("L2").Value = 320
$ws.Range("M2").Value = -438.375
$ws.Range("N2").Value = -546
# Row 97: H97,I97,K97,M97
$ws.Range("H97").Value = 879.7143
$ws.Range("I97").Value = 879.7143
$ws.Range("K97").Value = 879.7143
$ws.Range("M97").Value = -383.7143
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 8573
$ws.Range("I126").Value = 8688.833000000001
$ws.Range("K126").Value = 26066.499
$ws.Range("M126").Value = -23596.499

$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22,I22,J22,K22,L22,M22,N22
$ws.Range("H22").Value = 3833.3333
$ws.Range("I22").Value = 3833.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3833.3333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3538.3333
$ws.Range("N22").ClearContents()
# Row 27: H27,I27,J27,K27,L27,M27,N27
$ws.Range("H27").Value = 3833.3333
$ws.Range("I27").Value = 3833.3333
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3833.3333
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3726.3333
$ws.Range("N27").ClearContents()
# Row 55: H55,I55,J55,K55,L55,M55,N55
$ws.Range("H55").Value = 347.42856
$ws.Range("I55").Value = 331.66666
$ws.Range("J55").Value = 359.25
$ws.Range("K55").Value = 331.66666
$ws.Range("L55").Value = 359.25
$ws.Range("M55").Value = -158.66666
$ws.Range("N55").Value = -705.25
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 3317.5789
$ws.Range("I136").Value = 2919.75
$ws.Range("K136").Value = 8759.25
$ws.Range("M136").Value = -6209.25

$ws = $wb.Worksheets.Item("WVR")
# Row 15: H15,J15,L15,N15
$ws.Range("H15").Value = 2511247.8
$ws.Range("J15").Value = 14997.333
$ws.Range("L15").Value = 14997.333
$ws.Range("N15").Value = -15573.333
# Row 62: H62,J62,L62,N62
$ws.Range("H62").Value = 4847.5713
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
# Row 65: H65,J65,L65,N65
$ws.Range("H65").Value = 4847.5713
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740
# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 1029.0625
$ws.Range("I107").Value = 685.7692
$ws.Range("K107").Value = 2057.3076
$ws.Range("M107").Value = -137.3076000000001
# Row 113: H113,I113,K113,M113
$ws.Range("H113").Value = 999.5
$ws.Range("I113").Value = 999.5
$ws.Range("K113").Value = 2998.5
$ws.Range("M113").Value = -828.5
# Row 124: H124,J124,L124,N124
$ws.Range("H124").Value = 32999.668
$ws.Range("J124").Value = 32999.668
$ws.Range("L124").Value = 32999.668
$ws.Range("N124").Value = -42819.668
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 1998.5
$ws.Range("I132").Value = 1998.5
$ws.Range("K132").Value = 5995.5
$ws.Range("M132").Value = -3465.5
# Row 135: H135,I135,K135,M135
$ws.Range("H135").Value = 750000
$ws.Range("I135").Value = 750000
$ws.Range("K135").Value = 750000
$ws.Range("M135").Value = -744930

